$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.576.56'
$ws.Range('E2').Value = '  +8.17%  '
$ws.Range('D3').Value = '1.766.95'
$ws.Range('E3').Value = '  +4.44%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '223.60'
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('D6').Value = '0.549'
$ws.Range('E6').Value = '  +3.09%  '
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '30.20'
$ws.Range('E8').Value = '  +2.52%  '
$ws.Range('D9').Value = '46.57'
$ws.Range('E9').Value = '  +3.94%  '
$ws.Range('D10').Value = '0.275'
$ws.Range('E10').Value = '  +3.28%  '
$ws.Range('D11').Value = '0.0654'
$ws.Range('E11').Value = '  +2.25%  '
$ws.Range('D12').Value = '0.0922'
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('D13').Value = '2.028.94'
$ws.Range('E13').Value = '  +4.84%  '
$ws.Range('D14').Value = '1.776.60'
$ws.Range('E14').Value = '  +4.88%  '
$ws.Range('D15').Value = '0.622'
$ws.Range('E15').Value = '  +2.18%  '
$ws.Range('D16').Value = '33.622.46'
$ws.Range('E16').Value = '  +8.29%  '
$ws.Range('D17').Value = '9.88'
$ws.Range('E17').Value = '  -4.03%  '
$ws.Range('D18').Value = '4.15'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = '67.92'
$ws.Range('E19').Value = '  +1.53%  '
$ws.Range('D20').Value = '249.37'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = '0.0₃0732'
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').Value = '10.17'
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('D24').Value = '4.16'
$ws.Range('D25').Value = '2.13'
$ws.Range('E25').Value = '  -1.34%  '
$ws.Range('D26').Value = '158.14'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = '16.32'
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('D28').Value = '0.113'
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('D29').Value = '6.85'
$ws.Range('E29').Value = '  +1.76%  '
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('D31').Value = '3.76'
$ws.Range('E31').Value = '  +5.40%  '
$ws.Range('D32').Value = '0.0508'
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('E33').Value = '  +2.93%  '
$ws.Range('D34').Value = '3.51'
$ws.Range('E34').Value = '  +4.78%  '
$ws.Range('D35').Value = '1.471.75'
$ws.Range('E35').Value = '  -2.96%  '
$ws.Range('D36').Value = '1.77'
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('E37').Value = '  +2.68%  '
$ws.Range('D38').Value = '0.627'
$ws.Range('E38').Value = '  +1.85%  '
$ws.Range('D39').Value = '82.66'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('D40').Value = '0.0183'
$ws.Range('E40').Value = '  +1.98%  '
$ws.Range('E41').Value = '  +2.62%  '
$ws.Range('D42').Value = '2.68'
$ws.Range('E42').Value = '  +0.51%  '
$ws.Range('D43').Value = '0.876'
$ws.Range('E43').Value = '  +3.57%  '
$ws.Range('D44').Value = '2.05'
$ws.Range('E44').Value = '  +0.95%  '
$ws.Range('D45').Value = '0.0507'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('E46').Value = '  +3.21%  '
$ws.Range('D47').Value = '1.928.41'
$ws.Range('E47').Value = '  +5.53%  '
$ws.Range('D48').Value = '5.73'
$ws.Range('E48').Value = '  +2.84%  '
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').Value = '11.82'
$ws.Range('E50').Value = '  +13.68%  '
$ws.Range('D51').Value = '50.21'
$ws.Range('E51').Value = '  -3.21%  '
